$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    "black yoga leggings",
    "black yoga leggings with pockets for women",
    "black yoga pant",
    "black yoga pants long",
    "black yoga pants women capri",
    "black yoga pants women high waist",
    "black yoga tights for women",
    "blacks leggings",
    "blacks tights",
    "bleach body wash",
    "bleach skin for dark area",
    "blends chart",
    "blessed be leggings",
    "blue compresion pants",
    "blue compression leggings for women",
    "blue compression shorts women",
    "blue compression tights for boys",
    "blue line wedding rings for men",
    "blue running pants women",
    "blue running tights",
    "blue shorts running",
    "blue thigh high tights",
    "blue womens tights",
    "blue womens trousers",
    "blues super rugby",
    "body by victoria pants",
    "body chaffing stick",
    "body dynamic",
    "body dynamics workout",
    "body leggings for women",
    "body muscle chart",
    "body of a spartan",
    "body pants for women",
    "body rolls tight pants",
    "body spartan",
    "body spartan apparel",
    "body spartan clothing men",
    "body sport back brace",
    "body sport band",
    "body sport pilates ring",
    "body tape for skin",
    "body tape for women",
    "body tights",
    "body tights for women",
    "body workout",
    "body wrap it works",
    "bomb squad game",
    "bombs bullets and bribes",
    "bone and joint support",
    "bone tights for women",
    "bonito compression",
    "born primitive",
    "bottom tights",
    "boxing ankle support",
    "boxing ankle weights",
    "boxing apperal",
    "boxing ball on string",
    "boxing braces",
    "boxing compression",
    "boxing equipment for women",
    "boxing final",
    "boxing fitness tracker",
    "boxing frame",
    "boxing gear for women",
    "boxing gear women",
    "boxing gym shirt",
    "boxing pants women",
    "boxing set for women",
    "boxing shies",
    "boxing shorts for women",
    "boxing shorts with pockets",
    "boxing shorts women",
    "boxing target ball",
    "boxing tennis ball",
    "boxing torso",
    "boxing track suit",
    "boxing training shorts",
    "boxing wear",
    "boy compression leggings",
    "boy compression pants basketball",
    "boy flight jacket",
    "boys basketball compression tights",
    "boys basketball leggings",
    "boys blue compression tights",
    "boys comoression pants",
    "boys compession pants",
    "boys compression leggings",
    "boys compression pants",
    "boys compression pants black",
    "boys compression running pants",
    "boys compression running shorts",
    "boys compression tight",
    "boys compression tights basketball",
    "boys flight jacket",
    "boys full length compression leggings",
    "boys knee length shorts",
    "boys leggings",
    "boys running pants size 7",
    "boys running tank",
    "boys running tights"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

